$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "20.617.90"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.482.98"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9704"
$ws.Range("E5").Value = "  +2.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "279.72"
$ws.Range("E6").Value = "  -0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3662"
$ws.Range("E7").Value = "  -1.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3096"
$ws.Range("E8").Value = "  -2.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "40.17"
$ws.Range("E9").Value = "  -2.98%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.066"
$ws.Range("E10").Value = "  +0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06697"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("E12").Value = "  +0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.541"
$ws.Range("E13").Value = "  -0.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.15"
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.231"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  +2.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001033"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.480.87"
$ws.Range("E18").Value = "  +0.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.05963"
$ws.Range("E19").Value = "  +3.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "70.01"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.521"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.57"
$ws.Range("E22").Value = "  -0.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.07"
$ws.Range("E23").Value = "  -1.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.263"
$ws.Range("E24").Value = "  -0.31%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "20.666.36"
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.29"
$ws.Range("E26").Value = "  +3.32%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.123"
$ws.Range("E27").Value = "  -7.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.35"
$ws.Range("E28").Value = "  -1.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.639.41"
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "114.77"
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.937"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.058"
$ws.Range("E32").Value = "  -4.76%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.8269"
$ws.Range("E33").Value = "  -1.70%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08012"
$ws.Range("E34").Value = "  +2.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.552"
$ws.Range("E35").Value = "  -4.38%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.219"
$ws.Range("E36").Value = "  +8.26%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05808"
$ws.Range("E37").Value = "  -3.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.768"
$ws.Range("E38").Value = "  -3.22%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02050"
$ws.Range("E39").Value = "  -0.82%  "
$ws.Range("B40").Value = "Frax"
$ws.Range("C40").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9702"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.683"
$ws.Range("E41").Value = "  +3.08%  "
$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.47"
$ws.Range("E42").Value = "  -2.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1886"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5319"
$ws.Range("E44").Value = "  -1.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.540"
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.25"
$ws.Range("E46").Value = "  -1.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "118.69"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5216"
$ws.Range("E48").Value = "  -1.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.820"
$ws.Range("E49").Value = "  -0.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06517"
$ws.Range("E50").Value = "  +0.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.9865"
$ws.Range("E51").Value = "  -0.49%  "
